$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados..." timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 20:00"

# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 6189186
$ws.Range("C4").Value = 15950
$ws.Range("D4").Value = 3429680
$ws.Range("E4").Value = 2572082
$ws.Range("G4").Value = 200
$ws.Range("H4").Value = 187424

# Row 6 (India)
$ws.Range("B6").Value = 3684339
$ws.Range("C6").Value = 65170
$ws.Range("D6").Value = 2834973
$ws.Range("E6").Value = 783939
$ws.Range("G6").Value = 810
$ws.Range("H6").Value = 65427

# Row 20 (Francia)
$ws.Range("B20").Value = 281025
$ws.Range("C20").Value = 3082
$ws.Range("E20").Value = 164213
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 30635

# Row 22 (Italia)
$ws.Range("B22").Value = 269214
$ws.Range("C22").Value = 996
$ws.Range("D22").Value = 207653
$ws.Range("E22").Value = 26078
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 35483

# Row 23 (Alemania)
$ws.Range("B23").Value = 244566
$ws.Range("C23").Value = 1271
$ws.Range("E23").Value = 17354

# Row 27 (Canada)
$ws.Range("B27").Value = 128194
$ws.Range("C27").Value = 254
$ws.Range("D27").Value = 113790
$ws.Range("E27").Value = 5284
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 9120

# Row 32 (Ecuador)
$ws.Range("B32").Value = 113767
$ws.Range("C32").Value = 119
$ws.Range("D32").Value = 101669
$ws.Range("E32").Value = 5542
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 6556

# Row 35 (Republica Dominicana)
$ws.Range("B35").Value = 94715
$ws.Range("C35").Value = 474
$ws.Range("D35").Value = 68082
$ws.Range("E35").Value = 24923
$ws.Range("G35").Value = 29
$ws.Range("H35").Value = 1710

# Row 49 (Marruecos)
$ws.Range("B49").Value = 62590
$ws.Range("C49").Value = 1191
$ws.Range("D49").Value = 47595
$ws.Range("E49").Value = 13854
$ws.Range("G49").Value = 30
$ws.Range("H49").Value = 1141

# Row 54 (Barein)
$ws.Range("B54").Value = 52131
$ws.Range("C54").Value = 1009
$ws.Range("D54").Value = 18994
$ws.Range("E54").Value = 32328
$ws.Range("G54").Value = 16
$ws.Range("H54").Value = 809

# Row 55 (Etiopia)
$ws.Range("B55").Value = 51574
$ws.Range("D55").Value = 48654
$ws.Range("E55").Value = 2730
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 190

# Row 57 (Ghana)
$ws.Range("B57").Value = 44494
$ws.Range("C57").Value = 348
$ws.Range("D57").Value = 31244
$ws.Range("E57").Value = 11740
$ws.Range("G57").Value = 9
$ws.Range("H57").Value = 1510

# Row 58 (Argelia)
$ws.Range("B58").Value = 44298
$ws.Range("C58").Value = 93
$ws.Range("D58").Value = 42963
$ws.Range("E58").Value = 1059
$ws.Range("H58").Value = 276

# Row 70 (Irlanda)
$ws.Range("B70").Value = 28811
$ws.Range("C70").Value = 51
$ws.Range("E70").Value = 3670

# Row 80 (Paraguay)
$ws.Range("B80").Value = 17308
$ws.Range("C80").Value = 438
$ws.Range("D80").Value = 4811
$ws.Range("E80").Value = 12330
$ws.Range("G80").Value = 7
$ws.Range("H80").Value = 167

# Row 81 (Dinamarca)
$ws.Range("B81").Value = 17105
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 9146
$ws.Range("E81").Value = 7651
$ws.Range("H81").Value = 308

# Row 82 (Libano)
$ws.Range("B82").Value = 16985
$ws.Range("C82").Value = 94
$ws.Range("D82").Value = 15205
$ws.Range("E82").Value = 1156
$ws.Range("H82").Value = 624

# Row 94 (Albania)
$ws.Range("B94").Value = 9513
$ws.Range("C94").Value = 133
$ws.Range("D94").Value = 5214
$ws.Range("E94").Value = 4015
$ws.Range("G94").Value = 4
$ws.Range("H94").Value = 284

# Row 102 (Maldivas)
$ws.Range("B102").Value = 7804
$ws.Range("C102").Value = 137
$ws.Range("D102").Value = 5155
$ws.Range("E102").Value = 2621

# Row 114 (Suazilandia)
$ws.Range("B114").Value = 4577
$ws.Range("C114").Value = 16
$ws.Range("D114").Value = 3529
$ws.Range("E114").Value = 957

# Row 128 (Gambia)
$ws.Range("B128").Value = 2972
$ws.Range("C128").Value = 44
$ws.Range("D128").Value = 1288
$ws.Range("E128").Value = 1652
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = 32

# Row 129 (Uganda)
$ws.Range("B129").Value = 2963
$ws.Range("D129").Value = 1032
$ws.Range("E129").Value = 1835
$ws.Range("H129").Value = 96

# Row 132 (Mali)
$ws.Range("B132").Value = 2776
$ws.Range("C132").Value = 3
$ws.Range("E132").Value = 481

# Row 148 (Trinidad yTobago)
$ws.Range("B148").Value = 1724
$ws.Range("C148").Value = 91
$ws.Range("D148").Value = 493
$ws.Range("E148").Value = 1225
$ws.Range("H148").Value = 6

# Row 149 (Reunion)
$ws.Range("B149").Value = 1683
$ws.Range("D149").Value = 672
$ws.Range("E149").Value = 990
$ws.Range("H149").Value = 21

# Row 150 (Botsuana)
$ws.Range("B150").Value = 1634
$ws.Range("D150").Value = 880
$ws.Range("E150").Value = 745
$ws.Range("H150").Value = 9

# Row 164 (Republica del Chad)
$ws.Range("B164").Value = 1013
$ws.Range("C164").Value = 1
$ws.Range("D164").Value = 880
$ws.Range("E164").Value = 56
